# Update the cryptos list with the latest scraped prices / volume deltas.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$c = $ws.Cells.Item(2, 4)
$c.Value = "'60.837.69"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  -0.22%  "

# Row 3 - Ethereum
$c = $ws.Cells.Item(3, 4)
$c.Value = "'2.910.51"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -0.09%  "

# Row 4 - TetherUSD
$ws.Cells.Item(4, 5).Value = "  +0.07%  "

# Row 5 - BNB
$c = $ws.Cells.Item(5, 4)
$c.Value = "'591.90"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.95%  "

# Row 6 - Solana
$c = $ws.Cells.Item(6, 4)
$c.Value = "'145.50"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.71%  "

# Row 7 - USDC
$ws.Cells.Item(7, 5).Value = "  +0.06%  "

# Row 8 - XRP
$ws.Cells.Item(8, 5).Value = "  +0.65%  "

# Row 9 - Toncoin
$c = $ws.Cells.Item(9, 4)
$c.Value = "'6.87"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +1.35%  "

# Row 10 - Dogecoin
$ws.Cells.Item(10, 5).Value = "  -0.98%  "

# Row 11 - Cardano
$c = $ws.Cells.Item(11, 4)
$c.Value = "'0.439"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -2.17%  "

# Row 12 - ShibaInu
$ws.Cells.Item(12, 5).Value = "  +0.08%  "

# Row 13 - Avalanche
$c = $ws.Cells.Item(13, 4)
$c.Value = "'33.47"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -0.56%  "

# Row 14 - TRON
$c = $ws.Cells.Item(14, 4)
$c.Value = "'0.126"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -0.63%  "

# Row 15 - WrappedliquidstakedEther2.0
$c = $ws.Cells.Item(15, 4)
$c.Value = "'3.393.29"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -0.09%  "

# Row 16 - WrappedBTC
$c = $ws.Cells.Item(16, 4)
$c.Value = "'60.872.12"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -0.03%  "

# Row 17 - Polkadot
$c = $ws.Cells.Item(17, 4)
$c.Value = "'6.67"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -1.42%  "

# Row 18 - WrappedEther
$c = $ws.Cells.Item(18, 4)
$c.Value = "'2.913.56"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +0.05%  "

# Row 19 - BitcoinCash
$c = $ws.Cells.Item(19, 4)
$c.Value = "'429.34"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.44%  "

# Row 20 - Chainlink
$c = $ws.Cells.Item(20, 4)
$c.Value = "'13.33"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -1.98%  "

# Row 21 - Polygon
$c = $ws.Cells.Item(21, 4)
$c.Value = "'0.677"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +0.68%  "

# Row 22 - Uniswap
$c = $ws.Cells.Item(22, 4)
$c.Value = "'7.03"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -1.28%  "

# Row 23 - Litecoin
$c = $ws.Cells.Item(23, 4)
$c.Value = "'81.41"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +1.37%  "

# Row 24 - RenderToken
$ws.Cells.Item(24, 5).Value = "  -0.33%  "

# Row 25 - Fetch.AI
$c = $ws.Cells.Item(25, 4)
$c.Value = "'2.19"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.85%  "

# Row 26 - InternetComputer(DFINITY)
$c = $ws.Cells.Item(26, 4)
$c.Value = "'11.85"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -0.24%  "

# Row 27 - Dai
$ws.Cells.Item(27, 5).Value = "  -0.01%  "

# Row 28 - ImmutableX
$c = $ws.Cells.Item(28, 4)
$c.Value = "'2.27"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +5.10%  "

# Row 29 - FirstDigitalUSD
$ws.Cells.Item(29, 5).Value = "  +0.10%  "

# Row 30 - PancakeSwap
$ws.Cells.Item(30, 5).Value = "  -0.59%  "

# Row 31 - NEARProtocol
$c = $ws.Cells.Item(31, 4)
$c.Value = "'7.03"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -3.02%  "

# Row 32 - EthereumClassic
$c = $ws.Cells.Item(32, 4)
$c.Value = "'26.45"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -0.06%  "

# Row 33 - Hedera
$ws.Cells.Item(33, 5).Value = "  +0.65%  "

# Row 34 - PEPE
$c = $ws.Cells.Item(34, 4)
$c.Value = "'0.0₃0849"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +0.35%  "

# Row 35 - Mantle
$ws.Cells.Item(35, 5).Value = "  -0.08%  "

# Row 36 - Filecoin
$c = $ws.Cells.Item(36, 4)
$c.Value = "'5.61"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -0.43%  "

# Row 37 - dogwifhat
$ws.Cells.Item(37, 5).Value = "  +0.46%  "

# Row 38 - Kaspa
$ws.Cells.Item(38, 5).Value = "  -1.17%  "

# Row 39 - Stacks
$c = $ws.Cells.Item(39, 4)
$c.Value = "'1.97"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -1.89%  "

# Row 40 - Cosmos
$c = $ws.Cells.Item(40, 4)
$c.Value = "'8.51"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -1.82%  "

# Row 41 - TheGraph
$c = $ws.Cells.Item(41, 4)
$c.Value = "'0.285"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -2.02%  "

# Row 42 - Arweave
$c = $ws.Cells.Item(42, 4)
$c.Value = "'39.89"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -4.24%  "

# Row 43 - Bittensor
$c = $ws.Cells.Item(43, 4)
$c.Value = "'373.51"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -0.98%  "

# Row 46 - Monero
$c = $ws.Cells.Item(46, 4)
$c.Value = "'132.35"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -0.57%  "

# Row 47 - USDe
$ws.Cells.Item(47, 5).Value = "  -0.11%  "

# Row 48 - InjectiveProtocol
$c = $ws.Cells.Item(48, 4)
$c.Value = "'23.73"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -4.27%  "

# Row 49 - Stellar
$ws.Cells.Item(49, 5).Value = "  -0.58%  "

# Row 50 - ThetaToken
$c = $ws.Cells.Item(50, 4)
$c.Value = "'2.00"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -3.82%  "

# Row 51 - Cronos
$ws.Cells.Item(51, 5).Value = "  +0.87%  "

# Row 44 - was Maker, now VeChain (rows 44/45 swapped position in ranking)
$ws.Cells.Item(44, 2).Value = "VeChain"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Cells.Item(44, 4)
$c.Value = "'0.0343"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -1.07%  "

# Row 45 - was VeChain, now Maker
$ws.Cells.Item(45, 2).Value = "Maker"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$c = $ws.Cells.Item(45, 4)
$c.Value = "'2.697.66"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +0.96%  "
